$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D3 uses its own (non-shared) copy of the PASS/FAIL formula - clarify the
# status report so a B-column error shows up as "ERROR" instead of being
# silently compared (and reported FAIL) against C.
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'

# D4:D41 share a single formula definition (si="0") rooted at D4; updating
# D4's formula text re-derives the same logic (relative refs) for every
# dependent cell in the shared group, including the error rows (D16, D22:D27,
# D29:D31) whose B-column precedent is itself an error value - those now
# recompute to the string "ERROR" rather than propagating the #NUM! error.
$ws.Range("D4:D41").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'
